# Apply the edits described by the commit:
#  - Update a handful of text labels across the SERVER and CAM sheets.
#  - Shrink SERVER row 4's height to match its shorter wrapped text.
#  - Move the active selection in each sheet, and make SERVER the active tab.

$wb = $excel.ActiveWorkbook

$wsSite = $wb.Worksheets.Item("Siteground")
$wsServer = $wb.Worksheets.Item("SERVER")
$wsCam = $wb.Worksheets.Item("CAM")

# --- Text edits -----------------------------------------------------------

# SERVER sheet: "3. .htaccess..." -> "2. .htaccess..."
$wsServer.Range("A4").Value = "2. .htaccess Protected Folders & Files "

# SERVER sheet: "Ugly popup" -> " Popup" (first line of the Cons cell)
$wsServer.Range("H4").Value = " Popup" + [char]10 + "No operator/supervisor logic"

# CAM sheet: drop " (Built-in)" suffix
$wsCam.Range("A4").Value = "3. Camera BasicAuth"

# CAM sheet: drop " (BEST)" suffix
$wsCam.Range("A5").Value = "6. Server-Hosted Camera UI"

# --- Row height -------------------------------------------------------------

# SERVER row 4 text got shorter, so its wrapped height shrinks.
$wsServer.Rows.Item(4).RowHeight = 37.75

# --- Selections / active sheet ---------------------------------------------

# Non-active sheets just remember their own last selection.
$null = $wsSite.Range("L16").Select()
$null = $wsCam.Range("A5").Select()

# SERVER becomes the active sheet with A4 selected.
$null = $wsServer.Activate()
$null = $wsServer.Range("A4").Select()
